$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Convert E83:E85 from text-stored numbers to real numeric values.
$ws.Range("E83").Value = 20
$ws.Range("E84").Value = 533274
$ws.Range("E85").Value = 526371

# Append three new breakout rows (86-88) from the latest screener run.
# E column keeps bsecode as text (as the source JSON/CSV feed emits it) -
# a leading apostrophe forces text storage; re-applying the Normal style
# afterwards drops the auto-added "quote prefix" number format so the
# cell's style index stays untouched, matching the other plain data cells.
$ws.Range("A86").Value = "27/06/2024 09:44:43"
$ws.Range("B86").Value = 1
$ws.Range("C86").Value = "PRESTIGE"
$ws.Range("D86").Value = "Prestige Estates Projects Limited"
$ws.Range("E86").Value = "'533274"
$ws.Range("E86").Style = "Normal"
$ws.Range("F86").Value = -3.98
$ws.Range("G86").Value = 1850.45
$ws.Range("H86").Value = 1173409

$ws.Range("A87").Value = "27/06/2024 09:44:43"
$ws.Range("B87").Value = 2
$ws.Range("C87").Value = "OBEROIRLTY"
$ws.Range("D87").Value = "Oberoi Realty Limited"
$ws.Range("E87").Value = "'533273"
$ws.Range("E87").Style = "Normal"
$ws.Range("F87").Value = -2.2
$ws.Range("G87").Value = 1756.9
$ws.Range("H87").Value = 1164911

$ws.Range("A88").Value = "27/06/2024 09:44:43"
$ws.Range("B88").Value = 3
$ws.Range("C88").Value = "NMDC"
$ws.Range("D88").Value = "Nmdc Limited"
$ws.Range("E88").Value = "'526371"
$ws.Range("E88").Style = "Normal"
$ws.Range("F88").Value = -1.56
$ws.Range("G88").Value = 245.35
$ws.Range("H88").Value = 15717847
